# Add VPC figures (#12)
# Insert a new "Sensitivity" section (sensXls / sensSheet) into the
# Workflow sheet, just above the existing "TaskdoSensitivityAnalysis"
# section, by inserting 3 new rows at row 22 and filling them in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Make room: insert 3 blank rows at row 22 (old rows 22-26 shift to 25-29).
#    EntireRow.Insert() on a 3-row range inserts all three rows in one go,
#    and the new rows inherit the format of row 21 (A:s9 / B:s5).
$ws.Range("A22:A24").EntireRow.Insert()

# 2) Row 22: new section header "Sensitivity" (same look as the other
#    section header rows, e.g. row 3 "simulation").
$ws.Range("A3:F3").Copy()
$ws.Range("A22:F22").PasteSpecial(-4122)
$ws.Range("B22").Value = "Sensitivity"
$ws.Rows.Item(22).RowHeight = 49.8

# 3) Row 23: sensXls
$ws.Range("A23").Value = "sensXls"
$ws.Range("B23").Value = "xlsfilefor sensitivity Parameter definition; if it is empty, sheet is in this xlsfile"
$ws.Range("D18:F18").Copy()
$ws.Range("D23:F23").PasteSpecial(-4122)
$ws.Rows.Item(23).RowHeight = 26.4

# 4) Row 24: sensSheet
$ws.Range("A24").Value = "sensSheet"
$ws.Range("B24").Value = "xlssheet for sensitivity Parameter definition; if empty first sheet is taken"
$ws.Range("D18:F18").Copy()
$ws.Range("D24:F24").PasteSpecial(-4122)
$ws.Rows.Item(24).RowHeight = 26.4

# 5) Row 25 (old row 22, the existing "TaskdoSensitivityAnalysis" header,
#    shifted down by the insert) keeps its content/style but its height
#    changes from 49.8 to 13.8.
$ws.Rows.Item(25).RowHeight = 13.8

# 6) Row 28 (old row 25, shifted down by the insert) no longer needs its
#    explicit 13.8 custom height - the sheet's default row height is
#    already 13.8, so drop back to the auto/default height.
$ws.Rows.Item(28).AutoFit()

# 7) Restore the selection to the cell shown in the final workbook.
$ws.Range("C24").Select()
